# VG-FS-ADO-Sync.xlsx : "ADO->FS ID field key" fix
#
# The "SingleField" sheet mapped an erroneous FS field "source_control_reference"
# to ADO field "id" (direction ADO_TO_FS) in row 9. That row is bogus / a
# duplicate of the real "source_control_reference_created_on" -> "System.CreatedDate"
# mapping that already exists one row below it, so it is removed outright
# (the rows below shift up).
#
# Two ADO-Field-Key values were also missing their "Custom." prefix, and a
# stray space in "System. State" is corrected to "System.State".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SingleField")

# Fix ADO-Field-Key values that were missing the "Custom." prefix.
$ws.Range("C4").Value = "Custom.IMSCreatedOn"
$ws.Range("C5").Value = "Custom.IMSQueuedOn"

# Remove the bogus row: source_control_reference | true | id | ADO_TO_FS
$ws.Rows(9).Delete()

# Fix typo "System. State" -> "System.State" (now on row 10 after the delete,
# in column C, the ADO-Field-Key for "devops_status").
$ws.Range("C10").Value = "System.State"

# Keep the active selection in sync with the now-shorter used range
# (previously highlighted rows 9:10, now just row 9).
[void]$ws.Range("A9:XFD9").Select()
